# Fix bug: not removing ID and TIME in features. Rerun training and model assessment.
# Append two new result rows (row 7 and row 8) to the metrics sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for the two additional model runs
$newRows = @(
    @{ Index = 5; Model = "LogisticRegression";     Acc = 0.5949599465954606; Prec = 0.5416560980135188; Rec = 0.5949599465954606; F1 = 0.5457186484909325; Hamming = 0.4050400534045394; LogLoss = 0.972586098896736;  Roc = 0.5029389521611491 },
    @{ Index = 6; Model = "RandomForestClassifier"; Acc = 0.5797174009790832; Prec = 0.5351378178425146; Rec = 0.5797174009790832; F1 = 0.5435588035343082; Hamming = 0.4202825990209168; LogLoss = 0.7260962108894713; Roc = 0.4979390002818844 }
)

$startRow = 7

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($r, 1).Value2 = $data.Index
    $ws.Cells.Item($r, 2).Value2 = $data.Model
    $ws.Cells.Item($r, 3).Value2 = $data.Acc
    $ws.Cells.Item($r, 4).Value2 = $data.Prec
    $ws.Cells.Item($r, 5).Value2 = $data.Rec
    $ws.Cells.Item($r, 6).Value2 = $data.F1
    $ws.Cells.Item($r, 7).Value2 = $data.Hamming
    $ws.Cells.Item($r, 8).Value2 = $data.LogLoss
    $ws.Cells.Item($r, 9).Value2 = $data.Roc

    # Column A uses the bold/bordered/centered style used by the other index cells (e.g. A2)
    $ws.Cells.Item(2, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
